$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" (the last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# --- Header row (bold, centered, thin border - matches the other sheets) ---
$ws.Cells.Item(1,1).Value = "ds"
$ws.Cells.Item(1,2).Value = "PO_Forecast"
$ws.Cells.Item(1,3).Value = "yhat_lower"
$ws.Cells.Item(1,4).Value = "yhat_upper"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows ---
$data = @(
    @(45025.99999999999, 70, -57.60856723560526, 197.12242250094),
    @(45032.99999999999, 71, -48.81965012255916, 194.331159964443),
    @(45046.99999999999, 73, -55.86446618759714, 209.6828617504705),
    @(45053.99999999999, 74, -51.09990030024967, 202.6499034232732),
    @(45067.99999999999, 76, -54.03589468046058, 205.2570916728295),
    @(45074.99999999999, 77, -53.13631300799709, 202.2113623324555),
    @(45081.99999999999, 78, -52.51307762231691, 193.4479367638591),
    @(45088.99999999999, 79, -56.36301886260774, 207.154132257562),
    @(45095.99999999999, 80, -47.91654568790253, 211.2622146244834),
    @(45102.99999999999, 81, -40.80488660396878, 200.279716854368),
    @(45109.99999999999, 82, -47.47722395137555, 204.5246730149104),
    @(45123.99999999999, 84, -39.68313071784359, 207.1038983226011),
    @(45130.99999999999, 84, -34.22028902806881, 206.0014867182419),
    @(45137.99999999999, 85, -41.33744722882256, 203.864874519841),
    @(45144.99999999999, 86, -41.14589274586984, 215.8039896036346),
    @(45151.99999999999, 87, -49.54589927081017, 208.162260909397),
    @(45158.99999999999, 88, -36.67983870093261, 209.3379111970686),
    @(45165.99999999999, 89, -35.92492353466768, 223.3341457507771),
    @(45172.99999999999, 90, -30.10267239315218, 217.7974295866267),
    @(45179.99999999999, 91, -32.27250225922514, 218.7852674897666),
    @(45186.99999999999, 92, -39.32360290896273, 218.3050518163494),
    @(45193.99999999999, 93, -44.45623783864209, 216.7491256828592),
    @(45214.99999999999, 96, -28.80956400658151, 217.9556104405984),
    @(45221.99999999999, 97, -26.35990577684787, 222.0132826990824),
    @(45228.99999999999, 98, -35.29124029419165, 219.2018384201755),
    @(45235.99999999999, 99, -22.88068274467107, 225.7934827004731),
    @(45242.99999999999, 100, -30.66749780173316, 223.2409796052005),
    @(45249.99999999999, 101, -24.03583138942167, 235.9856403704042),
    @(45256.99999999999, 101, -29.00524281572673, 233.1820860187117),
    @(45263.99999999999, 102, -26.41873036276019, 235.4642775508224),
    @(45277.99999999999, 104, -19.24400464868381, 237.2922065011009),
    @(45298.99999999999, 107, -25.56879016145766, 235.9690933926576),
    @(45312.99999999999, 109, -22.25958876908686, 236.0139849805244),
    @(45354.99999999999, 115, -11.24412133765992, 240.5929232925409),
    @(45424.99999999999, 124, -1.341832401101204, 241.201024396021),
    @(45459.99999999999, 129, 9.89846468862263, 249.1591546020886),
    @(45466.99999999999, 130, 11.30976529547469, 258.210596501035),
    @(45487.99999999999, 133, 4.120943213235334, 257.4006392365955),
    @(45494.99999999999, 134, 11.66014310237502, 263.2643585454342),
    @(45501.99999999999, 135, 8.257796550263581, 259.6578292632533),
    @(45508.99999999999, 136, 6.055649324885635, 261.3496916582055),
    @(45515.99999999999, 136, 13.87335073449718, 263.7094316517279),
    @(45522.99999999999, 137, 6.183401511203953, 255.9384453550017),
    @(45529.99999999999, 138, 9.362262891394, 263.1858571722306),
    @(45536.99999999999, 139, 13.50064397415771, 262.7352625706256),
    @(45578.99999999999, 145, 18.57954366602796, 269.2577063493997),
    @(45585.99999999999, 146, 24.27161532585775, 276.078348669432),
    @(45592.99999999999, 147, 20.18703768672599, 270.8569423361853),
    @(45599.99999999999, 148, 12.8608120476648, 268.9626278092605),
    @(45606.99999999999, 149, 24.13246722295109, 280.5973856525353),
    @(45613.99999999999, 150, 32.65542475089024, 279.0693938689324),
    @(45620.99999999999, 151, 24.2653284678802, 276.0305747842062),
    @(45627.99999999999, 152, 28.10779013292497, 279.6672857079956),
    @(45634.99999999999, 153, 18.13048627926886, 276.6678841638175)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
}

# Date formatting for column A (matches the "Order Week"/"Order Month" date style)
$ws.Range("A2:A55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
